$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text cells to remain text (avoid Excel auto-converting numeric-looking
# strings like "243.34" into real numbers), then drop the temporary number
# format so the cell style reverts back to its original (unstyled) state.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '36.561.85'
$ws.Range('E2').Value = '  +0.70%  '
Set-TextValue $ws.Range('D3') '1.943.07'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('E4').Value = '  -0.34%  '
Set-TextValue $ws.Range('D5') '243.34'
$ws.Range('E5').Value = '  +0.40%  '
Set-TextValue $ws.Range('D6') '0.612'
$ws.Range('E6').Value = '  +0.46%  '
$ws.Range('E7').Value = '  -0.18%  '
Set-TextValue $ws.Range('D8') '57.21'
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('E9').Value = '  -0.79%  '
$ws.Range('E10').Value = '  -1.49%  '
$ws.Range('E11').Value = '  -0.61%  '
Set-TextValue $ws.Range('D12') '21.77'
$ws.Range('E12').Value = '  +2.59%  '
Set-TextValue $ws.Range('D13') '2.227.60'
$ws.Range('E13').Value = '  -0.24%  '
Set-TextValue $ws.Range('D14') '0.806'
$ws.Range('E14').Value = '  -2.25%  '
Set-TextValue $ws.Range('D15') '13.32'
$ws.Range('E15').Value = '  -1.09%  '
$ws.Range('E16').Value = '  -1.35%  '
Set-TextValue $ws.Range('D17') '1.952.29'
$ws.Range('E17').Value = '  +1.27%  '
Set-TextValue $ws.Range('D18') '36.454.93'
$ws.Range('E18').Value = '  +0.61%  '
Set-TextValue $ws.Range('D19') '69.28'
$ws.Range('E19').Value = '  -0.68%  '
Set-TextValue $ws.Range('D20') '0.0₃0855'
$ws.Range('E20').Value = '  -1.18%  '
Set-TextValue $ws.Range('D21') '227.53'
$ws.Range('E21').Value = '  -0.44%  '
Set-TextValue $ws.Range('D22') '4.96'
$ws.Range('E22').Value = '  -0.32%  '
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('E24').Value = '  -4.23%  '
Set-TextValue $ws.Range('D25') '2.29'
$ws.Range('E25').Value = '  +1.21%  '
Set-TextValue $ws.Range('D26') '9.19'
$ws.Range('E26').Value = '  -2.03%  '
Set-TextValue $ws.Range('D27') '159.49'
$ws.Range('E27').Value = '  -2.25%  '
$ws.Range('E28').Value = '  +14.97%  '
Set-TextValue $ws.Range('D29') '19.19'
$ws.Range('E29').Value = '  -0.67%  '
$ws.Range('E30').Value = '  +0.52%  '
$ws.Range('E31').Value = '  -4.25%  '
$ws.Range('E32').Value = '  -1.25%  '
Set-TextValue $ws.Range('D33') '0.0616'
$ws.Range('E33').Value = '  -2.04%  '
Set-TextValue $ws.Range('D34') '4.16'
$ws.Range('E34').Value = '  -2.77%  '
Set-TextValue $ws.Range('D35') '6.14'
$ws.Range('E35').Value = '  +1.75%  '
$ws.Range('E36').Value = '  -0.20%  '
$ws.Range('E37').Value = '  -0.90%  '
$ws.Range('E38').Value = '  +2.78%  '
Set-TextValue $ws.Range('D39') '3.28'
$ws.Range('E39').Value = '  +15.16%  '
$ws.Range('E40').Value = '  +1.44%  '
$ws.Range('E41').Value = '  +1.32%  '
$ws.Range('E42').Value = '  -0.31%  '
$ws.Range('E43').Value = '  -2.94%  '
$ws.Range('E44').Value = '  +1.33%  '
Set-TextValue $ws.Range('D45') '1.345.06'
$ws.Range('E45').Value = '  +0.58%  '
Set-TextValue $ws.Range('D46') '1.03'
$ws.Range('E46').Value = '  -0.25%  '
Set-TextValue $ws.Range('D47') '86.25'
$ws.Range('E47').Value = '  -1.61%  '
$ws.Range('E48').Value = '  -3.18%  '
$ws.Range('E49').Value = '  +0.24%  '
Set-TextValue $ws.Range('D50') '2.120.54'
$ws.Range('E50').Value = '  -0.29%  '
Set-TextValue $ws.Range('D51') '43.15'
$ws.Range('E51').Value = '  -6.53%  '
